# Fruta / hortaliza, semanal
#
# A new weekly price record (two quality grades: "Maduro" and "Pintón") for
# "Femacal de La Calera" / Plátano is inserted right after the existing
# row 600. This pushes the previously-existing rows 601-650 down to
# 603-652 and grows the sheet's used range from A1:T650 to A1:T652.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 601 (shifts old rows 601..650 -> 603..652)
$ws.Rows("601:602").Insert()

# --- New row 601: Plátano, Maduro ---
$ws.Range("A601").Value = 3
$ws.Range("B601").Value = "Femacal de La Calera"
$ws.Range("C601").Value = "Coquimbo"
$ws.Range("D601").Value = 44578
$ws.Range("E601").Value = 5
$ws.Range("F601").Value = "Fruta"
$ws.Range("G601").Value = 100108
$ws.Range("H601").Value = "Tropicales y subtropicales"
$ws.Range("I601").Value = 100108006
$ws.Range("J601").Value = "Plátano"
$ws.Range("K601").Value = "Sin especificar"
$ws.Range("L601").Value = "Maduro"
$ws.Range("M601").Value = 130
$ws.Range("N601").Value = 12000
$ws.Range("O601").Value = 12000
$ws.Range("P601").Value = 12000
$ws.Range("Q601").Value = "$/caja 20 kilos"
$ws.Range("R601").Value = "Ecuador"
$ws.Range("S601").Value = 600
$ws.Range("T601").Value = 20

# --- New row 602: Plátano, Pintón ---
$ws.Range("A602").Value = 3
$ws.Range("B602").Value = "Femacal de La Calera"
$ws.Range("C602").Value = "Coquimbo"
$ws.Range("D602").Value = 44578
$ws.Range("E602").Value = 5
$ws.Range("F602").Value = "Fruta"
$ws.Range("G602").Value = 100108
$ws.Range("H602").Value = "Tropicales y subtropicales"
$ws.Range("I602").Value = 100108006
$ws.Range("J602").Value = "Plátano"
$ws.Range("K602").Value = "Sin especificar"
$ws.Range("L602").Value = "Pintón"
$ws.Range("M602").Value = 250
$ws.Range("N602").Value = 14000
$ws.Range("O602").Value = 15000
$ws.Range("P602").Value = 14480
$ws.Range("Q602").Value = "$/caja 20 kilos"
$ws.Range("R602").Value = "Ecuador"
$ws.Range("S602").Value = 724
$ws.Range("T602").Value = 20

# Keep column D's date number format consistent with the rest of the column.
$ws.Range("D601:D602").NumberFormat = $ws.Range("D600").NumberFormat
